$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PlainText($addr, $val) {
    $c = $ws.Range($addr)
    $c.Value = "'" + $val
    $c.Style = "Normal"
}

$ws.Range("D2").Value = "64.471.22"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "3.152.34"
$ws.Range("E3").Value = "  -0.96%  "
$ws.Range("E4").Value = "  -0.12%  "
Set-PlainText "D5" "612.02"
$ws.Range("E5").Value = "  +0.66%  "
Set-PlainText "D6" "144.66"
$ws.Range("E6").Value = "  -1.06%  "
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("D8").Value = "3.149.73"
$ws.Range("E8").Value = "  -0.01%  "
Set-PlainText "D9" "0.529"
$ws.Range("E9").Value = "  +0.56%  "
Set-PlainText "D10" "0.152"
$ws.Range("E10").Value = "  +0.58%  "
Set-PlainText "D11" "5.40"
$ws.Range("E11").Value = "  -1.96%  "
Set-PlainText "D12" "0.474"
$ws.Range("E12").Value = "  +0.27%  "
Set-PlainText "D13" "0.0000258"
$ws.Range("E13").Value = "  +1.53%  "
Set-PlainText "D14" "35.64"
$ws.Range("E14").Value = "  -0.79%  "
$ws.Range("D15").Value = "3.670.27"
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("E16").Value = "  +2.95%  "
$ws.Range("D17").Value = "64.440.28"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").Value = "3.152.54"
$ws.Range("E18").Value = "  -0.35%  "
Set-PlainText "D19" "6.88"
$ws.Range("E19").Value = "  -0.68%  "
Set-PlainText "D20" "478.18"
$ws.Range("E20").Value = "  +0.12%  "
Set-PlainText "D21" "14.72"
$ws.Range("E21").Value = "  +1.22%  "
Set-PlainText "D22" "0.728"
$ws.Range("E22").Value = "  +2.15%  "
Set-PlainText "D23" "7.89"
$ws.Range("E23").Value = "  +2.27%  "
Set-PlainText "D24" "13.76"
$ws.Range("E24").Value = "  +0.30%  "
Set-PlainText "D25" "84.86"
$ws.Range("E25").Value = "  +1.55%  "
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-PlainText "D27" "8.68"
$ws.Range("E27").Value = "  +3.76%  "
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-PlainText "D28" "2.81"
$ws.Range("E28").Value = "  -2.67%  "
Set-PlainText "D29" "7.51"
$ws.Range("E29").Value = "  +10.99%  "
Set-PlainText "D30" "0.119"
$ws.Range("E30").Value = "  +0.06%  "
Set-PlainText "D31" "2.11"
$ws.Range("E31").Value = "  -3.48%  "
$ws.Range("E32").Value = "  -0.01%  "
Set-PlainText "D33" "26.72"
$ws.Range("E33").Value = "  +2.19%  "
Set-PlainText "D34" "2.67"
$ws.Range("E34").Value = "  -2.88%  "
$ws.Range("E35").Value = "  +1.54%  "
Set-PlainText "D36" "5.98"
$ws.Range("E36").Value = "  -0.22%  "
$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").Value = "0.0₃0753"
$ws.Range("E37").Value = "  +5.18%  "
$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-PlainText "D38" "52.76"
$ws.Range("E38").Value = "  -2.68%  "
Set-PlainText "D39" "3.10"
$ws.Range("E39").Value = "  +6.82%  "
Set-PlainText "D40" "456.95"
$ws.Range("E40").Value = "  +1.25%  "
Set-PlainText "D41" "0.0398"
$ws.Range("E41").Value = "  +0.84%  "
$ws.Range("E42").Value = "  +0.64%  "
Set-PlainText "D43" "8.37"
$ws.Range("E43").Value = "  -0.59%  "
$ws.Range("D44").Value = "2.866.47"
$ws.Range("E44").Value = "  +0.59%  "
Set-PlainText "D45" "0.270"
$ws.Range("E45").Value = "  +1.08%  "
Set-PlainText "D46" "2.28"
$ws.Range("E46").Value = "  +1.59%  "
Set-PlainText "D47" "2.45"
$ws.Range("E47").Value = "  +6.54%  "
Set-PlainText "D48" "26.61"
$ws.Range("E48").Value = "  +1.24%  "
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("E50").Value = "  +0.62%  "
Set-PlainText "D51" "33.90"
$ws.Range("E51").Value = "  +5.00%  "
